$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 1000
$ws.Range("K12").Value = 1000
$ws.Range("M12").Value = -830
$ws.Range("H21").Value = 1420
$ws.Range("I21").Value = 1420
$ws.Range("K21").Value = 1420
$ws.Range("M21").Value = -952
$ws.Range("H23").Value = 1420
$ws.Range("I23").Value = 1420
$ws.Range("K23").Value = 1420
$ws.Range("M23").Value = -1186
$ws.Range("H29").Value = 1237.5555
$ws.Range("I29").Value = 84.5
$ws.Range("J29").Value = 2160
$ws.Range("K29").Value = 253.5
$ws.Range("L29").Value = 6480
$ws.Range("M29").Value = 27.5
$ws.Range("N29").Value = -7042
$ws.Range("H38").Value = 328.6111
$ws.Range("I38").Value = 89.61539
$ws.Range("J38").Value = 950
$ws.Range("K38").Value = 268.84617
$ws.Range("L38").Value = 2850
$ws.Range("M38").Value = 103.15383
$ws.Range("N38").Value = -3594
$ws.Range("H46").Value = 9563.333000000001
$ws.Range("J46").Value = 9563.333000000001
$ws.Range("L46").Value = 28689.999
$ws.Range("N46").Value = -28927.999
$ws.Range("H58").Value = 5837.7407
$ws.Range("I58").Value = 284.26666
$ws.Range("J58").Value = 12779.583
$ws.Range("K58").Value = 852.79998
$ws.Range("L58").Value = 38338.749
$ws.Range("M58").Value = -702.79998
$ws.Range("N58").Value = -38638.749
$ws.Range("H60").Value = 9563.333000000001
$ws.Range("J60").Value = 9563.333000000001
$ws.Range("L60").Value = 28689.999
$ws.Range("N60").Value = -29657.999
$ws.Range("H87").Value = 35000
$ws.Range("J87").Value = 35000
$ws.Range("L87").Value = 35000
$ws.Range("N87").Value = -37496
$ws.Range("H90").Value = 35000
$ws.Range("J90").Value = 35000
$ws.Range("L90").Value = 105000
$ws.Range("N90").Value = -117480
$ws.Range("H107").Value = 459.5
$ws.Range("I107").Value = 432.5
$ws.Range("J107").Value = 479.75
$ws.Range("K107").Value = 432.5
$ws.Range("L107").Value = 479.75
$ws.Range("M107").Value = 1487.5
$ws.Range("N107").Value = -4319.75
$ws.Range("H116").Value = 1795.75
$ws.Range("I116").Value = 1397.9166
$ws.Range("J116").Value = 2094.125
$ws.Range("K116").Value = 1397.9166
$ws.Range("L116").Value = 2094.125
$ws.Range("M116").Value = 2044.0834
$ws.Range("N116").Value = -8978.125
$ws.Range("H132").Value = 2978927.2
$ws.Range("I132").Value = 3403925.8
$ws.Range("J132").Value = 3937.6667
$ws.Range("K132").Value = 10211777.4
$ws.Range("L132").Value = 11813.0001
$ws.Range("M132").Value = -10209247.4
$ws.Range("N132").Value = -16873.0001
$ws.Range("H137").Value = 1233.9565
$ws.Range("I137").Value = 1114.7894
$ws.Range("J137").Value = 1800
$ws.Range("K137").Value = 3344.3682
$ws.Range("L137").Value = 5400
$ws.Range("M137").Value = -794.3681999999999
$ws.Range("N137").Value = -10500
$ws.Range("H141").Value = 2224.4443
$ws.Range("I141").Value = 1615.3334
$ws.Range("K141").Value = 4846.0002
$ws.Range("M141").Value = 333.9997999999996
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 418.2
$ws.Range("I74").Value = 389.6154
$ws.Range("J74").Value = 604
$ws.Range("K74").Value = 389.6154
$ws.Range("L74").Value = 604
$ws.Range("M74").Value = 484.3846
$ws.Range("N74").Value = -2352
$ws.Range("H77").Value = 418.2
$ws.Range("I77").Value = 389.6154
$ws.Range("J77").Value = 604
$ws.Range("K77").Value = 1948.077
$ws.Range("L77").Value = 3020
$ws.Range("M77").Value = 2419.923
$ws.Range("N77").Value = -11756
$ws.Range("H110").Value = 4546.9395
$ws.Range("I110").Value = 7704.9414
$ws.Range("J110").Value = 1191.5625
$ws.Range("K110").Value = 7704.9414
$ws.Range("L110").Value = 1191.5625
$ws.Range("M110").Value = -5659.9414
$ws.Range("N110").Value = -5281.5625
$ws.Range("H132").Value = 3382.1428
$ws.Range("I132").Value = 3340.1191
$ws.Range("J132").Value = 3466.1904
$ws.Range("K132").Value = 10020.3573
$ws.Range("L132").Value = 10398.5712
$ws.Range("M132").Value = -7490.3573
$ws.Range("N132").Value = -15458.5712
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2099.3333
$ws.Range("I105").Value = 2030.7142
$ws.Range("J105").Value = 2195.4
$ws.Range("K105").Value = 2030.7142
$ws.Range("L105").Value = 2195.4
$ws.Range("M105").Value = -283.7141999999999
$ws.Range("N105").Value = -5689.4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2101.8333
$ws.Range("I16").Value = 1152.75
$ws.Range("K16").Value = 1152.75
$ws.Range("M16").Value = -865.75
$ws.Range("H99").Value = 2179.08
$ws.Range("J99").Value = 2482.4285
$ws.Range("L99").Value = 2482.4285
$ws.Range("N99").Value = -5478.4285
$ws.Range("H113").Value = 2101.8333
$ws.Range("I113").Value = 1152.75
$ws.Range("K113").Value = 1152.75
$ws.Range("M113").Value = 1017.25
$ws.Range("H126").Value = 2179.08
$ws.Range("J126").Value = 2482.4285
$ws.Range("L126").Value = 7447.2855
$ws.Range("N126").Value = -12387.2855
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2718.3333
$ws.Range("I3").Value = 2077.5
$ws.Range("J3").Value = 4000
$ws.Range("K3").Value = 6232.5
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = -6120.5
$ws.Range("N3").Value = -12224
$ws.Range("H5").Value = 599.9355
$ws.Range("I5").Value = 583.1579
$ws.Range("J5").Value = 626.5
$ws.Range("K5").Value = 1749.4737
$ws.Range("L5").Value = 1879.5
$ws.Range("M5").Value = -1637.4737
$ws.Range("N5").Value = -2103.5
$ws.Range("H34").Value = 1696
$ws.Range("I34").Value = 799.5
$ws.Range("J34").Value = 1994.8334
$ws.Range("K34").Value = 2398.5
$ws.Range("L34").Value = 5984.5002
$ws.Range("M34").Value = -2314.5
$ws.Range("N34").Value = -6152.5002
$ws.Range("H39").Value = 2990
$ws.Range("J39").Value = 2990
$ws.Range("L39").Value = 8970
$ws.Range("N39").Value = -9558
$ws.Range("H55").Value = 3349.75
$ws.Range("J55").Value = 3349.75
$ws.Range("L55").Value = 10049.25
$ws.Range("N55").Value = -10403.25
$ws.Range("H82").Value = 6351.6665
$ws.Range("I82").Value = 1166.6666
$ws.Range("J82").Value = 7092.381
$ws.Range("K82").Value = 3499.9998
$ws.Range("L82").Value = 21277.143
$ws.Range("M82").Value = -3093.9998
$ws.Range("N82").Value = -22089.143
$ws.Range("H85").Value = 6351.6665
$ws.Range("I85").Value = 1166.6666
$ws.Range("J85").Value = 7092.381
$ws.Range("K85").Value = 3499.9998
$ws.Range("L85").Value = 21277.143
$ws.Range("M85").Value = -2095.9998
$ws.Range("N85").Value = -24085.143
$ws.Range("H107").Value = 773.86365
$ws.Range("I107").Value = 550
$ws.Range("J107").Value = 901.7857
$ws.Range("K107").Value = 1650
$ws.Range("L107").Value = 2705.3571
$ws.Range("M107").Value = 270
$ws.Range("N107").Value = -6545.3571
$ws.Range("H126").Value = 3323.077
$ws.Range("J126").Value = 3745.4546
$ws.Range("L126").Value = 11236.3638
$ws.Range("N126").Value = -21116.3638
$ws.Range("H131").Value = 4369596.5
$ws.Range("I131").Value = 9530
$ws.Range("J131").Value = 6985636
$ws.Range("K131").Value = 28590
$ws.Range("L131").Value = 20956908
$ws.Range("M131").Value = -23550
$ws.Range("N131").Value = -20966988
$ws.Range("H132").Value = 1834.5714
$ws.Range("J132").Value = 2408.889
$ws.Range("L132").Value = 21680.001
$ws.Range("N132").Value = -26740.001
$ws.Range("H135").Value = 599.9355
$ws.Range("I135").Value = 583.1579
$ws.Range("J135").Value = 626.5
$ws.Range("K135").Value = 5248.4211
$ws.Range("L135").Value = 5638.5
$ws.Range("M135").Value = -2713.4211
$ws.Range("N135").Value = -10708.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6241
$ws.Range("I122").Value = 6241
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 18723
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -16273
$ws.Range("N122").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2741.0833
$ws.Range("I7").Value = 2399.6
$ws.Range("J7").Value = 2985
$ws.Range("K7").Value = 2399.6
$ws.Range("L7").Value = 2985
$ws.Range("M7").Value = -2287.6
$ws.Range("N7").Value = -3209
$ws.Range("H126").Value = 2741.0833
$ws.Range("I126").Value = 2399.6
$ws.Range("J126").Value = 2985
$ws.Range("K126").Value = 7198.799999999999
$ws.Range("L126").Value = 8955
$ws.Range("M126").Value = -4728.799999999999
$ws.Range("N126").Value = -13895
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 361428.56
$ws.Range("H81").Value = 2000
$ws.Range("I81").Value = 2000
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 4000
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -2939
$ws.Range("N81").Value = -6122
$ws.Range("H84").Value = 2000
$ws.Range("I84").Value = 2000
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 20000
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -14696
$ws.Range("N84").Value = -30608
$ws.Range("H113").Value = 719.875
$ws.Range("I113").Value = 859.8333
$ws.Range("J113").Value = 300
$ws.Range("K113").Value = 2579.4999
$ws.Range("L113").Value = 900
$ws.Range("M113").Value = -409.4998999999998
$ws.Range("N113").Value = -5240
$ws.Range("H132").Value = 1129.2046
$ws.Range("I132").Value = 832.2
$ws.Range("J132").Value = 4099.25
$ws.Range("K132").Value = 2496.6
$ws.Range("L132").Value = 12297.75
$ws.Range("M132").Value = 33.39999999999964
$ws.Range("N132").Value = -17357.75
